$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-24 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-25 Monday", 2) | Out-Null
$d.Content.Find.Execute("31+31=", $true, $false, $false, $false, $false, $true, 1, $false, "39+44=", 2) | Out-Null
$d.Content.Find.Execute("3+36=", $true, $false, $false, $false, $false, $true, 1, $false, "64-34=", 2) | Out-Null
$d.Content.Find.Execute("1+23=", $true, $false, $false, $false, $false, $true, 1, $false, "42+18=", 2) | Out-Null
$d.Content.Find.Execute("35-20=", $true, $false, $false, $false, $false, $true, 1, $false, "57+4=", 2) | Out-Null
$d.Content.Find.Execute("86-32=", $true, $false, $false, $false, $false, $true, 1, $false, "30-4=", 2) | Out-Null
$d.Content.Find.Execute("92-51=", $true, $false, $false, $false, $false, $true, 1, $false, "51+31=", 2) | Out-Null
$d.Content.Find.Execute("46+0=", $true, $false, $false, $false, $false, $true, 1, $false, "8+44=", 2) | Out-Null
$d.Content.Find.Execute("54+42=", $true, $false, $false, $false, $false, $true, 1, $false, "93-70=", 2) | Out-Null
$d.Content.Find.Execute("99-72=", $true, $false, $false, $false, $false, $true, 1, $false, "87+12=", 2) | Out-Null
$d.Content.Find.Execute("24-13=", $true, $false, $false, $false, $false, $true, 1, $false, "61+32=", 2) | Out-Null
$d.Content.Find.Execute("17+54=", $true, $false, $false, $false, $false, $true, 1, $false, "8+11=", 2) | Out-Null
$d.Content.Find.Execute("52+23=", $true, $false, $false, $false, $false, $true, 1, $false, "31-0=", 2) | Out-Null
$d.Content.Find.Execute("25-11=", $true, $false, $false, $false, $false, $true, 1, $false, "99-86=", 2) | Out-Null
$d.Content.Find.Execute("97-39=", $true, $false, $false, $false, $false, $true, 1, $false, "71+5=", 2) | Out-Null
$d.Content.Find.Execute("21+46=", $true, $false, $false, $false, $false, $true, 1, $false, "8+29=", 2) | Out-Null
$d.Content.Find.Execute("27+26=", $true, $false, $false, $false, $false, $true, 1, $false, "15+30=", 2) | Out-Null
$d.Content.Find.Execute("76-52=", $true, $false, $false, $false, $false, $true, 1, $false, "73-29=", 2) | Out-Null
$d.Content.Find.Execute("66-62=", $true, $false, $false, $false, $false, $true, 1, $false, "6+5=", 2) | Out-Null
$d.Content.Find.Execute("65-57=", $true, $false, $false, $false, $false, $true, 1, $false, "96-57=", 2) | Out-Null
$d.Content.Find.Execute("57-47=", $true, $false, $false, $false, $false, $true, 1, $false, "17+35=", 2) | Out-Null
$d.Content.Find.Execute("3+29=", $true, $false, $false, $false, $false, $true, 1, $false, "25+31=", 2) | Out-Null
$d.Content.Find.Execute("23+74=", $true, $false, $false, $false, $false, $true, 1, $false, "11-9=", 2) | Out-Null
$d.Content.Find.Execute("86-27=", $true, $false, $false, $false, $false, $true, 1, $false, "55+16=", 2) | Out-Null
$d.Content.Find.Execute("76+1=", $true, $false, $false, $false, $false, $true, 1, $false, "33+37=", 2) | Out-Null
$d.Content.Find.Execute("11+38=", $true, $false, $false, $false, $false, $true, 1, $false, "81-26=", 2) | Out-Null
$d.Content.Find.Execute("29+11=", $true, $false, $false, $false, $false, $true, 1, $false, "49+22=", 2) | Out-Null
$d.Content.Find.Execute("38+61=", $true, $false, $false, $false, $false, $true, 1, $false, "39+16=", 2) | Out-Null
$d.Content.Find.Execute("40+4=", $true, $false, $false, $false, $false, $true, 1, $false, "4+45=", 2) | Out-Null
$d.Content.Find.Execute("25+0=", $true, $false, $false, $false, $false, $true, 1, $false, "0+47=", 2) | Out-Null
$d.Content.Find.Execute("4+25=", $true, $false, $false, $false, $false, $true, 1, $false, "18+60=", 2) | Out-Null
$d.Content.Find.Execute("32-28=", $true, $false, $false, $false, $false, $true, 1, $false, "25+58=", 2) | Out-Null
$d.Content.Find.Execute("64-8=", $true, $false, $false, $false, $false, $true, 1, $false, "38+15=", 2) | Out-Null
$d.Content.Find.Execute("48-14=", $true, $false, $false, $false, $false, $true, 1, $false, "18-14=", 2) | Out-Null
$d.Content.Find.Execute("21-16=", $true, $false, $false, $false, $false, $true, 1, $false, "30-15=", 2) | Out-Null
$d.Content.Find.Execute("15+69=", $true, $false, $false, $false, $false, $true, 1, $false, "59-41=", 2) | Out-Null
$d.Content.Find.Execute("68+23=", $true, $false, $false, $false, $false, $true, 1, $false, "94-4=", 2) | Out-Null
$d.Content.Find.Execute("86-77=", $true, $false, $false, $false, $false, $true, 1, $false, "57+31=", 2) | Out-Null
$d.Content.Find.Execute("14+10=", $true, $false, $false, $false, $false, $true, 1, $false, "92-14=", 2) | Out-Null
$d.Content.Find.Execute("23+52=", $true, $false, $false, $false, $false, $true, 1, $false, "34+19=", 2) | Out-Null
$d.Content.Find.Execute("72+4=", $true, $false, $false, $false, $false, $true, 1, $false, "44+37=", 2) | Out-Null
$d.Content.Find.Execute("53-6=", $true, $false, $false, $false, $false, $true, 1, $false, "51-39=", 2) | Out-Null
$d.Content.Find.Execute("71-5=", $true, $false, $false, $false, $false, $true, 1, $false, "3+7=", 2) | Out-Null
$d.Content.Find.Execute("65+15=", $true, $false, $false, $false, $false, $true, 1, $false, "72-21=", 2) | Out-Null
$d.Content.Find.Execute("68+27=", $true, $false, $false, $false, $false, $true, 1, $false, "90-34=", 2) | Out-Null
$d.Content.Find.Execute("60+13=", $true, $false, $false, $false, $false, $true, 1, $false, "60-8=", 2) | Out-Null
$d.Content.Find.Execute("45-15=", $true, $false, $false, $false, $false, $true, 1, $false, "35+59=", 2) | Out-Null
$d.Content.Find.Execute("11+24=", $true, $false, $false, $false, $false, $true, 1, $false, "96-80=", 2) | Out-Null
$d.Content.Find.Execute("75+23=", $true, $false, $false, $false, $false, $true, 1, $false, "84-54=", 2) | Out-Null
$d.Content.Find.Execute("55+8=", $true, $false, $false, $false, $false, $true, 1, $false, "83-1=", 2) | Out-Null
$d.Content.Find.Execute("90-75=", $true, $false, $false, $false, $false, $true, 1, $false, "81-69=", 2) | Out-Null
$d.Content.Find.Execute("10+61=", $true, $false, $false, $false, $false, $true, 1, $false, "73-9=", 2) | Out-Null
$d.Content.Find.Execute("50+7=", $true, $false, $false, $false, $false, $true, 1, $false, "56-18=", 2) | Out-Null
$d.Content.Find.Execute("16-3=", $true, $false, $false, $false, $false, $true, 1, $false, "21+62=", 2) | Out-Null
$d.Content.Find.Execute("74-32=", $true, $false, $false, $false, $false, $true, 1, $false, "2+55=", 2) | Out-Null
$d.Content.Find.Execute("58+30=", $true, $false, $false, $false, $false, $true, 1, $false, "68-30=", 2) | Out-Null
$d.Content.Find.Execute("99-37=", $true, $false, $false, $false, $false, $true, 1, $false, "72-26=", 2) | Out-Null
$d.Content.Find.Execute("62+31=", $true, $false, $false, $false, $false, $true, 1, $false, "52-49=", 2) | Out-Null
$d.Content.Find.Execute("14+68=", $true, $false, $false, $false, $false, $true, 1, $false, "34+41=", 2) | Out-Null
$d.Content.Find.Execute("40+48=", $true, $false, $false, $false, $false, $true, 1, $false, "46+51=", 2) | Out-Null
$d.Content.Find.Execute("62-16=", $true, $false, $false, $false, $false, $true, 1, $false, "22-15=", 2) | Out-Null
$d.Content.Find.Execute("11+15=", $true, $false, $false, $false, $false, $true, 1, $false, "83-9=", 2) | Out-Null
$d.Content.Find.Execute("84-14=", $true, $false, $false, $false, $false, $true, 1, $false, "15+37=", 2) | Out-Null
$d.Content.Find.Execute("9+33=", $true, $false, $false, $false, $false, $true, 1, $false, "48+14=", 2) | Out-Null
$d.Content.Find.Execute("4+9=", $true, $false, $false, $false, $false, $true, 1, $false, "75-26=", 2) | Out-Null
$d.Content.Find.Execute("23-17=", $true, $false, $false, $false, $false, $true, 1, $false, "96-55=", 2) | Out-Null
$d.Content.Find.Execute("11+74=", $true, $false, $false, $false, $false, $true, 1, $false, "71-30=", 2) | Out-Null
$d.Content.Find.Execute("49+12=", $true, $false, $false, $false, $false, $true, 1, $false, "43+49=", 2) | Out-Null
$d.Content.Find.Execute("20+1=", $true, $false, $false, $false, $false, $true, 1, $false, "14+72=", 2) | Out-Null
$d.Content.Find.Execute("70-9=", $true, $false, $false, $false, $false, $true, 1, $false, "67+7=", 2) | Out-Null
$d.Content.Find.Execute("45-36=", $true, $false, $false, $false, $false, $true, 1, $false, "83-3=", 2) | Out-Null
$d.Content.Find.Execute("99-28=", $true, $false, $false, $false, $false, $true, 1, $false, "11+19=", 2) | Out-Null
$d.Content.Find.Execute("94-88=", $true, $false, $false, $false, $false, $true, 1, $false, "33+11=", 2) | Out-Null
$d.Content.Find.Execute("73+19=", $true, $false, $false, $false, $false, $true, 1, $false, "77-25=", 2) | Out-Null
$d.Content.Find.Execute("87-75=", $true, $false, $false, $false, $false, $true, 1, $false, "8+90=", 2) | Out-Null
$d.Content.Find.Execute("13+51=", $true, $false, $false, $false, $false, $true, 1, $false, "11-10=", 2) | Out-Null
$d.Content.Find.Execute("78-71=", $true, $false, $false, $false, $false, $true, 1, $false, "0+26=", 2) | Out-Null
$d.Content.Find.Execute("86-18=", $true, $false, $false, $false, $false, $true, 1, $false, "15+83=", 2) | Out-Null
$d.Content.Find.Execute("88-67=", $true, $false, $false, $false, $false, $true, 1, $false, "87-57=", 2) | Out-Null
$d.Content.Find.Execute("8+39=", $true, $false, $false, $false, $false, $true, 1, $false, "97-23=", 2) | Out-Null
$d.Content.Find.Execute("96-62=", $true, $false, $false, $false, $false, $true, 1, $false, "88-43=", 2) | Out-Null
$d.Content.Find.Execute("13+30=", $true, $false, $false, $false, $false, $true, 1, $false, "38-23=", 2) | Out-Null
$d.Content.Find.Execute("14+79=", $true, $false, $false, $false, $false, $true, 1, $false, "80-61=", 2) | Out-Null
$d.Content.Find.Execute("19+8=", $true, $false, $false, $false, $false, $true, 1, $false, "6+64=", 2) | Out-Null
$d.Content.Find.Execute("64+4=", $true, $false, $false, $false, $false, $true, 1, $false, "3+21=", 2) | Out-Null
$d.Content.Find.Execute("11+50=", $true, $false, $false, $false, $false, $true, 1, $false, "52-31=", 2) | Out-Null
$d.Content.Find.Execute("36+12=", $true, $false, $false, $false, $false, $true, 1, $false, "67+11=", 2) | Out-Null
$d.Content.Find.Execute("35+31=", $true, $false, $false, $false, $false, $true, 1, $false, "93-8=", 2) | Out-Null
$d.Content.Find.Execute("34+30=", $true, $false, $false, $false, $false, $true, 1, $false, "46-10=", 2) | Out-Null
$d.Content.Find.Execute("53+16=", $true, $false, $false, $false, $false, $true, 1, $false, "35-21=", 2) | Out-Null
$d.Content.Find.Execute("43+39=", $true, $false, $false, $false, $false, $true, 1, $false, "71-59=", 2) | Out-Null
$d.Content.Find.Execute("62+21=", $true, $false, $false, $false, $false, $true, 1, $false, "45+16=", 2) | Out-Null
$d.Content.Find.Execute("45+7=", $true, $false, $false, $false, $false, $true, 1, $false, "93-71=", 2) | Out-Null
$d.Content.Find.Execute("82-76=", $true, $false, $false, $false, $false, $true, 1, $false, "72-1=", 2) | Out-Null
$d.Content.Find.Execute("95-5=", $true, $false, $false, $false, $false, $true, 1, $false, "4+80=", 2) | Out-Null
$d.Content.Find.Execute("28-0=", $true, $false, $false, $false, $false, $true, 1, $false, "2+21=", 2) | Out-Null
$d.Content.Find.Execute("66-20=", $true, $false, $false, $false, $false, $true, 1, $false, "3+69=", 2) | Out-Null
$d.Content.Find.Execute("32+2=", $true, $false, $false, $false, $false, $true, 1, $false, "49+11=", 2) | Out-Null
$d.Content.Find.Execute("96-35=", $true, $false, $false, $false, $false, $true, 1, $false, "30+16=", 2) | Out-Null
$d.Content.Find.Execute("18+10=", $true, $false, $false, $false, $false, $true, 1, $false, "28+17=", 2) | Out-Null
$d.Content.Find.Execute("15+35=", $true, $false, $false, $false, $false, $true, 1, $false, "9+26=", 2) | Out-Null
